$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 228, shifting the existing rows 228:268 down to 229:269.
$ws.Rows.Item(228).Insert()

# Populate the newly-inserted row 228 with a new data record (same market/product
# as the surrounding rows, new date + price figures).
$ws.Range("A228").Value = 4
$ws.Range("B228").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C228").Value = "Los Lagos"
$ws.Range("D228").Value = 44504
$ws.Range("E228").Value = 10
$ws.Range("F228").Value = "Fruta"
$ws.Range("G228").Value = 100108
$ws.Range("H228").Value = "Tropicales y subtropicales"
$ws.Range("I228").Value = 100108006
$ws.Range("J228").Value = "Plátano"
$ws.Range("K228").Value = "Sin especificar"
$ws.Range("L228").Value = "Primera Pintón"
$ws.Range("M228").Value = 600
$ws.Range("N228").Value = 22000
$ws.Range("O228").Value = 23000
$ws.Range("P228").Value = 22500
$ws.Range("Q228").Value = "$/caja 20 kilos"
$ws.Range("R228").Value = "Ecuador"
$ws.Range("S228").Value = 1125
$ws.Range("T228").Value = 20
